$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text, even if it looks like a number
    # (e.g. "5.28", "0.999"), matching the inlineStr cells in the source data.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '63.208.66'
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").Value = '2.451.19'
$ws.Range("E3").Value = '  +0.54%  '
Set-TextValue $ws.Range("D4") '0.999'
$ws.Range("E4").Value = '  -0.17%  '
Set-TextValue $ws.Range("D5") '572.27'
$ws.Range("E5").Value = '  +0.90%  '
Set-TextValue $ws.Range("D6") '146.30'
$ws.Range("E6").Value = '  +0.41%  '
$ws.Range("E7").Value = '  +0.06%  '
Set-TextValue $ws.Range("D8") '0.537'
$ws.Range("E8").Value = '  +0.74%  '
$ws.Range("D9").Value = '2.448.01'
$ws.Range("E9").Value = '  +0.33%  '
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("E11").Value = '  +1.30%  '
Set-TextValue $ws.Range("D12") '5.28'
$ws.Range("E12").Value = '  -0.64%  '
Set-TextValue $ws.Range("D13") '0.356'
$ws.Range("E13").Value = '  +0.03%  '
Set-TextValue $ws.Range("D14") '27.01'
$ws.Range("E14").Value = '  +0.53%  '
Set-TextValue $ws.Range("D15") '0.0000179'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").Value = '2.892.88'
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '63.074.37'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").Value = '2.454.16'
$ws.Range("E18").Value = '  +1.11%  '
Set-TextValue $ws.Range("D19") '11.29'
$ws.Range("E19").Value = '  +0.18%  '
Set-TextValue $ws.Range("D20") '7.31'
$ws.Range("E20").Value = '  +4.49%  '
Set-TextValue $ws.Range("D21") '328.61'
$ws.Range("E21").Value = '  +1.28%  '
Set-TextValue $ws.Range("D22") '4.20'
$ws.Range("E22").Value = '  +0.64%  '
Set-TextValue $ws.Range("D23") '2.07'
$ws.Range("E23").Value = '  +12.88%  '
Set-TextValue $ws.Range("D24") '0.999'
Set-TextValue $ws.Range("D25") '65.77'
$ws.Range("E25").Value = '  -2.31%  '
Set-TextValue $ws.Range("D26") '616.22'
$ws.Range("E26").Value = '  +5.19%  '
Set-TextValue $ws.Range("D27") '8.96'
$ws.Range("E27").Value = '  +4.63%  '
Set-TextValue $ws.Range("D28") '0.0000102'
$ws.Range("E28").Value = '  +2.01%  '
$ws.Range("D29").Value = '2.565.12'
$ws.Range("E29").Value = '  +0.35%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D30") '1.50'
$ws.Range("E30").Value = '  +3.77%  '
$ws.Range("B31").Value = 'Binance-PegBSC-USD'
$ws.Range("C31").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D31") '1.00'
$ws.Range("E31").Value = '  +0.23%  '
Set-TextValue $ws.Range("D32") '8.26'
$ws.Range("E32").Value = '  -2.44%  '
Set-TextValue $ws.Range("D33") '0.142'
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("E34").Value = '  +0.81%  '
Set-TextValue $ws.Range("D35") '5.20'
$ws.Range("E35").Value = '  +6.99%  '
Set-TextValue $ws.Range("D36") '1.53'
$ws.Range("E36").Value = '  +1.49%  '
$ws.Range("E37").Value = '  +0.10%  '
Set-TextValue $ws.Range("D38") '0.380'
$ws.Range("E38").Value = '  -0.57%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Range("D39") '5.41'
$ws.Range("E39").Value = '  +0.82%  '
$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D40") '18.79'
$ws.Range("E40").Value = '  +0.12%  '
Set-TextValue $ws.Range("D41") '145.52'
$ws.Range("E41").Value = '  -1.73%  '
Set-TextValue $ws.Range("D42") '1.79'
$ws.Range("E42").Value = '  -2.49%  '
Set-TextValue $ws.Range("D43") '2.60'
$ws.Range("E43").Value = '  +6.52%  '
$ws.Range("E44").Value = '  -0.02%  '
Set-TextValue $ws.Range("D45") '41.88'
$ws.Range("E45").Value = '  +0.59%  '
Set-TextValue $ws.Range("D46") '148.62'
$ws.Range("E46").Value = '  -0.08%  '
Set-TextValue $ws.Range("D47") '3.77'
$ws.Range("E47").Value = '  +2.65%  '
Set-TextValue $ws.Range("D48") '21.19'
$ws.Range("E48").Value = '  +3.03%  '
Set-TextValue $ws.Range("D49") '0.0533'
$ws.Range("E49").Value = '  -0.45%  '
Set-TextValue $ws.Range("D50") '0.602'
$ws.Range("E50").Value = '  -0.04%  '
Set-TextValue $ws.Range("D51") '0.0232'
$ws.Range("E51").Value = '  +0.19%  '
